# Edit: "Add AUROC on CKD"
# Target shape: "Google Shape;96;p13" (7th shape on slide 1)
#   Paragraph (bullet) 2: "92% roc_auc in predicting Chronic Kidney Disease,1 year before onset"
#       -> "0.954% AUROC in predicting Chronic Kidney Disease,1 year before onset using CNN"
#   Paragraph (bullet) 3: "Used 8 ML & DL models on electronic health records from Insurance claims data"
#       -> "Used tree based models on EHR from claims data for extracting important features"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(7)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- 1) "92% " + "roc_auc" + " in " (3 runs) -> single run "0.954% AUROC in "
$f1 = $tr.Find("92% roc_auc in ")
$f1.Text = "0.954% AUROC in "

# --- 2) "1 year before onset" -> "1 year before onset using CNN"
$f2 = $tr.Find("1 year before onset")
$f2.Text = "1 year before onset using CNN"

# --- 3) "Used 8 ML & DL models " -> "Used tree based models "
$f3 = $tr.Find("Used 8 ML & DL models ")
$f3.Text = "Used tree based models "

# --- 4) "on electronic health records from Insurance claims data"
#        -> "on " + "EHR" + " from claims data for extracting important features"
$marker = $tr.Find("Used tree based models ")
$afterMarker = $marker.Start

$f4 = $tr.Find("on electronic health records from Insurance claims data", $afterMarker)
$f4.Text = "on EHR from claims data for extracting important features"

$rEHR = $tr.Find("EHR", $afterMarker)
$rEHR.Font.Bold = 0
$rEHR.Font.Italic = 0
$rEHR.Font.Underline = 0

Write-Host "Edits applied."
